$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string used in column B (referenced by rows 2-9) from "Remove Head" to "Add Event"
$ws.Range("B2:B9").Value = "Add Event"

# Delete rows 10, 11, 12 (Noise percentage 80, 90, 100) entirely
$ws.Range("A10:C12").EntireRow.Delete()
